# Update "想去人数" (want-to-go count) figures on two sheets that share
# the same underlying event rows ("展览" and "全部类型").

$wb = $excel.ActiveWorkbook

$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F2").Value = 260
$ws1.Range("F4").Value = 903
$ws1.Range("F5").Value = 536

$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F2").Value = 260
$ws4.Range("F4").Value = 903
$ws4.Range("F6").Value = 536
